$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Data entry: the row-7 task now records 1 "ingreso" in week 2 (column N).
# This feeds the chain of shared "difference" formulas across the row
# (O7, R7, U7, X7, AA7, AD7, AG7, AJ7, AM7, AP7, AS7, AV7, AY7 each
# recompute from 1 to 0), the running total (AZ7: 3 -> 4) and the
# remainder (BA7: 1 -> 0) -- all handled automatically by recalculation.
$ws.Range("N7").Value = 1

# Update the active-cell selection on the lower-right (frozen) pane to
# reflect where the user ended up after making the edit.
$ws.Range("O10").Select()

$wb.Save()
